$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('ADMINISTRACION DE BASES DE DATOS', 143.9074096679688),
    @('ADMINISTRACION DE REDES', 121.9426116943359),
    @('ADMINISTRACION DE SERVIDORES', 144.1614990234375),
    @('ALGORITMIA', 133.6002960205078),
    @('ALMACENES DE DATOS (DATA WAREHOUSE)', 151.3702545166016),
    @('BASES DE DATOS', 152.0135345458984),
    @('CLASIFICACION INTELIGENTE DE DATOS', 160.2016906738281),
    @('CONTROL DE PROYECTOS', 198.3759002685547),
    @('ESTADISTICA Y PROCESOS ESTOCASTICOS', 355.6754760742188),
    @('ESTRUCTURAS DE DATOS I', 38.69380187988281),
    @('ESTRUCTURAS DE DATOS II', 75.54373168945312),
    @('HIPERMEDIA', 195.336669921875),
    @('INGENIERIA DE SOFTWARE I', 39.55695724487305),
    @('INGENIERIA DE SOFTWARE II', 139.5456237792969),
    @('METODOS MATEMATICOS I', 323.9101257324219),
    @('METODOS MATEMATICOS II', 286.1632995605469),
    @('MINERIA DE DATOS', 160.5572967529297),
    @('PROGRAMACION', 269.4276428222656),
    @('PROGRAMACION PARA INTERNET', 358.6423034667969),
    @('SEGURIDAD DE LA INFORMACION', 196.5658569335938),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE ALGORITMIA', 270.6994018554688),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE BASES DE DATOS', 232.9805603027344),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS I', 35.28647613525391),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE ESTRUCTURAS DE DATOS II', 138.6548919677734),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE INGENIERIA DE SOFTWARE I', 281.482666015625),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS I', 314.3976745605469),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE METODOS MATEMATICOS II', 307.2843017578125),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE PROGRAMACION', 290.6789245605469),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE SISTEMAS BASADOS EN CONOCIMIENTO', 126.8496170043945),
    @('SEMINARIO DE SOLUCION DE PROBLEMAS DE USO, ADAPTACION, EXPLOTACION DE SISTEMAS OPERATIVOS', 177.7617034912109),
    @('SISTEMAS BASADOS EN CONOCIMIENTO', 134.867919921875),
    @('TEORIA DE LA COMPUTACION', 58.2849006652832),
    @('USO, ADAPTACION Y EXPLOTACION DE SISTEMAS OPERATIVOS', 160.5096130371094)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = [double]$data[$i][1]
}

# Remove the now-extra row 35 (previously 34 data rows + header, now 33 data rows + header)
$ws.Rows.Item(35).Delete()